$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 5 device-spec records (rows 2-6), each tagged
# lang_code = "eng". Every record needs an "fra" twin sitting directly above
# it. We keep each existing record in place (just retagging its lang_code to
# "fra") and insert a brand-new row directly *below* it to carry the "eng"
# copy of the same data -- inserting below (rather than above, which would
# butt a fresh row against the header's distinct border style and pick up a
# blended/incorrect style) keeps every data row on the same plain style as
# before.

$numRecords = 5
$numCols = 9
$xlPasteValues = -4163

# Snapshot the 5 existing records before touching anything. Value2 (not
# Value, which this host surfaces oddly for reads) gives correctly-typed
# scalars back.
$records = @()
for ($r = 2; $r -le ($numRecords + 1); $r++) {
    $row = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $row += , ($ws.Cells.Item($r, $c).Value2)
    }
    $records += , $row
}

# Insert one fresh row directly below each original record, bottom-up (row 6
# first, row 2 last) so a not-yet-processed row's number never shifts out
# from under us. Each original record therefore keeps its original row
# number (2,3,4,5,6 -> final 2,4,6,8,10) and the new partner row lands right
# after it (3,5,7,9,11).
for ($r = ($numRecords + 1); $r -ge 2; $r--) {
    $ws.Rows.Item($r + 1).Insert()
}

# Columns B (id) and G (min_driver_ver) are numeric but formatted as text
# ("@"); writing a number straight into .Value on such a cell silently
# stores it as text. Flip to General while writing, then restore "@" so the
# formatting matches the source exactly.
$numericCols = @(2, 7)

for ($i = 0; $i -lt $numRecords; $i++) {
    $fraRow = 2 + (2 * $i)
    $engRow = $fraRow + 1
    $record = $records[$i]

    # The "fra" row is the original row -- every cell except lang_code (A)
    # already holds the right data and style, so only A needs touching.
    $ws.Cells.Item($fraRow, 1).Value = "fra"

    # The "eng" row is brand new and empty; populate every column.
    for ($c = 1; $c -le $numCols; $c++) {
        $value = $record[$c - 1]
        $cell = $ws.Cells.Item($engRow, $c)

        if ($numericCols -contains $c) {
            $cell.NumberFormat = "General"
            $cell.Value = $value
            $cell.NumberFormat = "@"
        } elseif ($c -eq 9) {
            # Column I literal text is "TRUE". Assigning the bare string
            # "TRUE"/"FALSE" via .Value auto-converts the cell to a real
            # Boolean, so route it through a text formula + paste-values
            # instead, which keeps it as plain text without leaving a
            # formula behind or touching the cell's style.
            $cell.Formula = '=T("' + $value + '")'
            $cell.Copy()
            $cell.PasteSpecial($xlPasteValues)
        } else {
            $cell.Value = $value
        }
    }

    $ws.Cells.Item($engRow, 1).Value = "eng"
}

# The sheet no longer carries an autofilter (it covered A1:I6, now stale) --
# drop it and the hidden _FilterDatabase defined name that backs it (turning
# off AutoFilterMode alone removes the <autoFilter> element but leaves the
# workbook-level defined name behind).
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# Match the refreshed selection: I2:I11, active cell I2.
$ws.Range("I2:I11").Select()
